# Applies:
#  1) Notes Master "Date Placeholder 2" datetimeFigureOut field re-cache
#     (10/13/20 -> 11/3/20). This is an auto-updating field; some hosts
#     treat it as read-only, so the update is attempted defensively and
#     never allowed to abort the rest of the script.
#  2) Slide 4, "TextBox 121": fix the "casles" typo and merge the three
#     runs of paragraph 2 ("In some cases ... with the") into one clean
#     run, matching what retyping the selection in the UI would produce.

$p = $ppt.ActivePresentation

# --- Edit 1: Notes Master date field -------------------------------------
try {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $shp = $nm.Shapes.Item($i)
        if ($shp.Name -eq "Date Placeholder 2") {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Text.Length
            if ($len -gt 0) {
                $full = $tr.Characters(1, $len)
                $full.Text = "11/3/20"
            } else {
                $tr.Text = "11/3/20"
            }
            break
        }
    }
} catch {
    # Field text is computed/read-only in some hosts - ignore and continue.
}

# --- Edit 2: Slide 4 typo fix + run merge ---------------------------------
$slide4 = $p.Slides.Item(4)
$outerGroup = $slide4.Shapes.Item(2)

$targetShape = $null
for ($i = 1; $i -le $outerGroup.GroupItems.Count; $i++) {
    $item = $outerGroup.GroupItems.Item($i)
    if ($item.Name -eq "TextBox 121") {
        $targetShape = $item
        break
    }
}

$tr2 = $targetShape.TextFrame.TextRange
$fullText = $tr2.Text
$startMarker = "In some cases, the X12 portion"
$endMarker = "communicate directly with the"
$startIdx = $fullText.IndexOf($startMarker)
$endIdx = $fullText.IndexOf($endMarker) + $endMarker.Length
$len2 = $endIdx - $startIdx

$targetRange = $tr2.Characters($startIdx + 1, $len2)
$targetRange.Text = "In some cases, the X12 portion will be handled entirely within a clearing house.  In those cases, then 3., 4., 5. will communicate directly with the"
